$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, border, centered) from A1 onto the
# newly introduced header cells E1:H1 before filling in values.
$ws.Range("A1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$headers = @("file", "n_clusters", "init", "max_iter", "silhouette_score", "calinski_harabasz_score", "davies_bouldin_score", "parameters_combinations")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows: file, n_clusters, init, max_iter, silhouette, calinski, davies, params
$data = @(
    @("eval_datasets/test_1_eval.xlsx", 3, "random", 100, 0.3531438267773461, 9202.942408559138, 0.9403239801300892, '"3,random,100"'),
    @("eval_datasets/test_1_eval.xlsx", 3, "random", 200, 0.3531438267773461, 9202.942408559138, 0.9403239801300892, '"3,random,200"'),
    @("eval_datasets/test_1_eval.xlsx", 3, "random", 300, 0.3531438267773461, 9202.942408559138, 0.9403239801300892, '"3,random,300"'),
    @("eval_datasets/test_2_eval.xlsx", 3, "random", 100, 0.3191442676270201, 6821.110548040679, 0.9741397510366555, '"3,random,100"'),
    @("eval_datasets/test_2_eval.xlsx", 3, "random", 200, 0.3164361350539214, 6806.690709089491, 0.9785388671056104, '"3,random,200"'),
    @("eval_datasets/test_2_eval.xlsx", 3, "random", 300, 0.318762300729125, 6808.966785614843, 0.9747675584662162, '"3,random,300"'),
    @("eval_datasets/test_3_eval.xlsx", 3, "random", 100, 0.3473605365154773, 10057.06353694939, 0.9342918971710779, '"3,random,100"'),
    @("eval_datasets/test_3_eval.xlsx", 3, "random", 200, 0.3503787921370744, 10122.14601905551, 0.9316492199462204, '"3,random,200"'),
    @("eval_datasets/test_3_eval.xlsx", 3, "random", 300, 0.3473605365154773, 10057.06353694939, 0.9342918971710779, '"3,random,300"'),
    @("eval_datasets/test_4_eval.xlsx", 3, "random", 100, 0.2947414202622447, 7089.959420187719, 1.070300270675494, '"3,random,100"'),
    @("eval_datasets/test_4_eval.xlsx", 3, "random", 200, 0.3191175704960609, 6829.772616121577, 1.00158103781754, '"3,random,200"'),
    @("eval_datasets/test_4_eval.xlsx", 3, "random", 300, 0.3183030108675407, 6932.374615717033, 1.008816220664122, '"3,random,300"')
)

$rowIndex = 2
foreach ($row in $data) {
    for ($colIndex = 0; $colIndex -lt $row.Length; $colIndex++) {
        $ws.Cells.Item($rowIndex, $colIndex + 1).Value = $row[$colIndex]
    }
    $rowIndex++
}
